# Apply the "Fixed update to excel issue" edit:
#  1. Rename the "Requested quantity" header (B1) on both existing sheets
#     to the new metric-specific names.
#  2. Add a new "PO Forecast" worksheet (after "Monthly Trend") containing
#     the ds / PO_Forecast / yhat_lower / yhat_upper forecast table.

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# --- 1. Rename headers on the existing sheets -----------------------------
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet at the end -------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the look & feel (outline props, margins) of the other sheets.
$wsForecast.Outline.SummaryRow = 1
$wsForecast.Outline.SummaryColumn = 1

$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Reuse the bold/centered header style and the date-number style already
# present in the workbook (copied from the "Weekly Quantity" sheet) instead
# of creating brand-new style records.
$wsWeekly.Range("A1").Copy() | Out-Null
$wsForecast.Range("A1:D1").PasteSpecial(-4122) | Out-Null

$wsWeekly.Range("A2").Copy() | Out-Null
$wsForecast.Range("A2:A12").PasteSpecial(-4122) | Out-Null

$wsForecast.Application.CutCopyMode = $false

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Forecast data rows
$data = @(
    @(45354.99999999999, 0,  -16.55526626236709, 15.14207488368432),
    @(45375.99999999999, 19,   3.271306082132699, 35.77187365751858),
    @(45382.99999999999, 26,  10.49349507532911,  43.01004658720104),
    @(45389.99999999999, 33,  17.08313669104012,  49.60779582036829),
    @(45396.99999999999, 40,  23.67431028107071,  55.42018196944495),
    @(45403.99999999999, 47,  30.48056894515365,  63.08826148179114),
    @(45410.99999999999, 54,  38.44745364773441,  70.79243963274726),
    @(45417.99999999999, 61,  45.02149766263718,  76.11351999912303),
    @(45424.99999999999, 68,  52.79820476990073,  84.08989818457042),
    @(45431.99999999999, 75,  58.58925802627218,  92.30201936183217),
    @(45438.99999999999, 81,  65.00637726074039,  97.72692568273634)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Range("A$row").Value = $r[0]
    $wsForecast.Range("B$row").Value = $r[1]
    $wsForecast.Range("C$row").Value = $r[2]
    $wsForecast.Range("D$row").Value = $r[3]
    $row++
}

$wsWeekly.Select() | Out-Null
$wsWeekly.Range("A1").Select() | Out-Null
